# add default taxonomy based on config
#
# Inserts a new "Taxonomy" column before the existing "Country" column
# (old column F), shifting the old F:K columns to G:L, and fills the new
# column with the default taxonomy value "level 3" for every data row.
# Also normalizes the capitalization of a few "Category" values that moved
# into column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F; everything from F onward shifts right.
$ws.Columns.Item(6).EntireColumn.Insert()

# New column header + default values.
$ws.Range("F1").Value = "Taxonomy"
$ws.Range("F2:F11").Value = "level 3"

# Normalize capitalization of the category labels now living in column I.
$ws.Range("I2").Value = "Glass Sourcing"
$ws.Range("I3").Value = "Glass Sourcing"
$ws.Range("I4").Value = "Ingredient Sourcing"
$ws.Range("I5").Value = "Ingredient Sourcing"
$ws.Range("I6").Value = "Ingredient Sourcing"
$ws.Range("I7").Value = "Ingredient Sourcing"
$ws.Range("I10").Value = "Use And End Of Life"
$ws.Range("I11").Value = "Use And End Of Life"
